$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Type" column labels (A2:A6) and add new row 7 (A7)
$ws.Range("A2").Value = "Amarilla"
$ws.Range("A3").Value = "Carretera"
$ws.Range("A4").Value = "Montana"
$ws.Range("A5").Value = "Paseo"
$ws.Range("A6").Value = "Velo"
$ws.Range("A7").Value = "VTT"

# Update TotalCount (B), Sale (C), Discount (D), Profit (E) columns
$ws.Range("B2").Value = 94
$ws.Range("C2").Value = 17747116.06
$ws.Range("D2").Value = 1290163.4400000002
$ws.Range("E2").Value = 2814104.0600000005

$ws.Range("B3").Value = 93
$ws.Range("C3").Value = 13815307.885000002
$ws.Range("D3").Value = 1122212.6149999998
$ws.Range("E3").Value = 1826804.885

$ws.Range("B4").Value = 93
$ws.Range("C4").Value = 15390801.880000003
$ws.Range("D4").Value = 1159032.62
$ws.Range("E4").Value = 2114754.8799999994

$ws.Range("B5").Value = 202
$ws.Range("C5").Value = 33011143.95000001
$ws.Range("D5").Value = 2600518.0499999993
$ws.Range("E5").Value = 4797437.95

$ws.Range("B6").Value = 109
$ws.Range("C6").Value = 18250059.465
$ws.Range("D6").Value = 1576709.035
$ws.Range("E6").Value = 2305992.4650000003

$ws.Range("B7").Value = 109
$ws.Range("C7").Value = 20511921.019999996
$ws.Range("D7").Value = 1456612.4800000002
$ws.Range("E7").Value = 3034608.0200000005
